# Update relay parameter values on the "relays" sheet for rows 2 and 3,
# and move the active cell selection from G4 to H4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

# Row 2 updates
$ws.Range("G2").Value = 1200
$ws.Range("K2").Value = 0.5
$ws.Range("N2").Value = 0.8
$ws.Range("O2").Value = 1.2

# Row 3 updates
$ws.Range("G3").Value = 1200
$ws.Range("K3").Value = 0.5
$ws.Range("N3").Value = 0.8
$ws.Range("O3").Value = 1.2

# Update the active selection from G4 to H4
$ws.Activate()
$ws.Range("H4").Select()
